$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.094.42"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.562.39"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.94"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.67"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.560.61"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +3.01%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.91"
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.165.95"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.99"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.574.19"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.206.81"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.115"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("E19").Value = "  +7.72%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.26"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.02"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.703.55"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.10"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.60"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.557.92"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  -6.41%  "
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.87"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.60"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.10"
$ws.Range("E41").Value = "  +3.59%  "
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.73"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.52"
$ws.Range("E49").Value = "  +6.44%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -0.58%  "
